$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date value that was bumped by one day
# (45204 -> 45205, i.e. 2023-10-05 -> 2023-10-06) for every data row
# (rows 2 through 411) in an automatic daily update.
$ws.Range("C2:C411").Value = 45205
